$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell refs whose target "Price" text looks like a plain number (single
# decimal point, e.g. "588.95"). Excel auto-types such a literal as a number
# when assigned through .Value, which would both change the stored type and
# introduce floating point noise (e.g. 588.95 -> 588.95000000000005). The
# source file keeps these as literal text ("t=inlineStr"), so force a Text
# number format before the assignment and then restore the default style so
# the cells look-and-feel (and lack of an explicit style index) is unchanged.

$ws.Range("D2").Value = '62.469.85'
$ws.Range("E2").Value = '  +2.42%  '

$ws.Range("D3").Value = '2.933.80'
$ws.Range("E3").Value = '  +1.73%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '588.95'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.43%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.19'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.38%  '

$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("E8").Value = '  +3.23%  '

$ws.Range("D9").Value = '2.934.70'
$ws.Range("E9").Value = '  +1.79%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.13'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.73%  '

$ws.Range("E11").Value = '  +9.38%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.435'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.92%  '

$ws.Range("E13").Value = '  +7.74%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.29'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.35%  '

$ws.Range("E15").Value = '  -0.62%  '

$ws.Range("D16").Value = '3.425.73'
$ws.Range("E16").Value = '  +1.87%  '

$ws.Range("D17").Value = '62.489.98'
$ws.Range("E17").Value = '  +2.55%  '

$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '2.942.54'
$ws.Range("E18").Value = '  +1.06%  '

$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.63'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.43%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '432.77'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.84%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.47'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.84%  '

$ws.Range("E22").Value = '  +1.71%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.95'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.60%  '

$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.20'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.57%  '

$ws.Range("B25").Value = 'RenderToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.04'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.96%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.89'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.72%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.09'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.47%  '

$ws.Range("E28").Value = '  -0.04%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.24'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +9.73%  '

$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.15'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.57%  '

$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.57'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.82%  '

$ws.Range("E32").Value = '  +20.97%  '

$ws.Range("E33").Value = '  +5.76%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.08'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.00%  '

$ws.Range("E35").Value = '  +0.03%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.986'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.78%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.55'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.79%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.00%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '49.60'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.23%  '

$ws.Range("E40").Value = '  +5.77%  '

$ws.Range("E41").Value = '  +0.41%  '

$ws.Range("E42").Value = '  +0.74%  '

$ws.Range("E43").Value = '  +4.70%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '38.86'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.01%  '

$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '2.696.36'
$ws.Range("E45").Value = '  +1.28%  '

$ws.Range("B46").Value = 'Monero'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '135.25'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.67%  '

$ws.Range("E47").Value = '  +3.25%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '351.54'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.32%  '

$ws.Range("E50").Value = '  +2.48%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.43'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.99%  '
